$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format before writing, so numeric-looking
# price strings (e.g. "0.679", "241.18") are preserved verbatim as
# text instead of being coerced to numbers by Excel.
$dCol = $ws.Range("D2:D51")
$dCol.NumberFormat = "@"

$ws.Range('D2').Value = '44.213.97'
$ws.Range('E2').Value = '  +1.65%  '
$ws.Range('D3').Value = '2.354.84'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '0.679'
$ws.Range('E5').Value = '  +5.38%  '
$ws.Range('D6').Value = '241.18'
$ws.Range('E6').Value = '  +3.34%  '
$ws.Range('D7').Value = '74.29'
$ws.Range('E7').Value = '  +6.39%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '0.565'
$ws.Range('E9').Value = '  +22.88%  '
$ws.Range('E10').Value = '  +5.87%  '
$ws.Range('D11').Value = '31.55'
$ws.Range('E11').Value = '  +20.02%  '
$ws.Range('D12').Value = '7.37'
$ws.Range('E12').Value = '  +18.43%  '
$ws.Range('E13').Value = '  +2.80%  '
$ws.Range('D14').Value = '2.703.53'
$ws.Range('E14').Value = '  -0.99%  '
$ws.Range('E15').Value = '  +7.46%  '
$ws.Range('D16').Value = '0.909'
$ws.Range('E16').Value = '  +7.29%  '
$ws.Range('D17').Value = '2.356.79'
$ws.Range('E17').Value = '  -0.70%  '
$ws.Range('D18').Value = '44.415.27'
$ws.Range('E18').Value = '  +2.10%  '
$ws.Range('D19').Value = '0.0000102'
$ws.Range('E19').Value = '  +3.90%  '
$ws.Range('D20').Value = '6.68'
$ws.Range('E20').Value = '  +5.77%  '
$ws.Range('D21').Value = '77.96'
$ws.Range('E21').Value = '  +5.24%  '
$ws.Range('D22').Value = '255.57'
$ws.Range('E22').Value = '  +1.91%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').Value = '3.76'
$ws.Range('E24').Value = '  -5.28%  '
$ws.Range('E25').Value = '  +4.27%  '
$ws.Range('E26').Value = '  +7.03%  '
$ws.Range('E27').Value = '  +1.43%  '
$ws.Range('D28').Value = '22.54'
$ws.Range('E28').Value = '  -1.55%  '
$ws.Range('D29').Value = '174.82'
$ws.Range('E29').Value = '  +1.07%  '
$ws.Range('E30').Value = '  +2.67%  '
$ws.Range('E31').Value = '  +3.87%  '
$ws.Range('E32').Value = '  +4.70%  '
$ws.Range('D33').Value = '5.36'
$ws.Range('E33').Value = '  +8.09%  '
$ws.Range('D34').Value = '0.0752'
$ws.Range('E34').Value = '  +9.19%  '
$ws.Range('D35').Value = '5.33'
$ws.Range('E35').Value = '  +5.05%  '
$ws.Range('D36').Value = '3.88'
$ws.Range('E36').Value = '  +7.40%  '
$ws.Range('D37').Value = '2.44'
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('D38').Value = '6.55'
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('D39').Value = '0.0273'
$ws.Range('E39').Value = '  +7.66%  '
$ws.Range('D40').Value = '19.39'
$ws.Range('E40').Value = '  +5.09%  '
$ws.Range('D41').Value = '8.97'
$ws.Range('E41').Value = '  +0.39%  '
$ws.Range('E42').Value = '  +0.15%  '
$ws.Range('E43').Value = '  +3.61%  '
$ws.Range('D44').Value = '0.0999'
$ws.Range('E44').Value = '  +5.45%  '
$ws.Range('E45').Value = '  +13.58%  '
$ws.Range('D46').Value = '100.55'
$ws.Range('E46').Value = '  +1.90%  '
$ws.Range('D47').Value = '2.46'
$ws.Range('E47').Value = '  +10.01%  '
$ws.Range('E48').Value = '  -1.79%  '
$ws.Range('E49').Value = '  -0.88%  '
$ws.Range('D50').Value = '1.450.61'
$ws.Range('E50').Value = '  +0.09%  '
$ws.Range('B51').Value = 'TerraClassic'
$ws.Range('C51').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D51').Value = '0.000206'

# Restore the default (General) formatting on column D now that the
# values have been stored as text, so no stray number-format style
# is left behind on the cells.
$dCol.ClearFormats()
